$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Network")
$ws.Rows("7:8").Delete()
